$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.029969539563001, 1.2933141042731946, 0.13631759066352356, 3.434090680999769)"
$ws.Range("C2").Value = "JSU(-1.2474500170971483, 1.1484710786502483, 3.0598936084929917, 4.673528819582653)"
$ws.Range("D2").Value = "NIG(0.7441724452235371, 0.5502857929737022, 1.125292948650341, 2.869207629390771)"
$ws.Range("E2").Value = "NIG(1.3198970150049356, 1.0169418221737487, 3.7037682312695086, 6.000799517282518)"
